$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells: AC1=Wins, AD1=Losses, AE1=Ties
# Copy style from an existing header cell (AB1) so formatting matches (bold, border, centered)
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill in team record values for all data rows (2-48)
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 29).Value = 80   # AC = column 29
    $ws.Cells.Item($r, 30).Value = 82   # AD = column 30
    $ws.Cells.Item($r, 31).Value = 0    # AE = column 31
}
